# Weekly update: a new price record (week of 2021-11-15) is published for
# "Vega Modelo de Temuco - Níspero". The new record becomes the first data
# row (row 2) and every existing record is pushed down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 2; this shifts the existing
# rows 2-9 down to rows 3-10, which is exactly what the new data requires.
$ws.Rows.Item(2).Insert()

# The freshly inserted row inherited some formatting from the row above it;
# strip that so the new row starts out with the workbook's default (unstyled)
# cells, matching every other data row.
$ws.Range("A2:T2").ClearFormats()

# Column D holds dates, so give the new row's date cell the same date
# number format used by the rest of the "Fecha" column.
$ws.Range("D2").NumberFormat = $ws.Range("D3").NumberFormat

# Populate the new row with the latest weekly record.
$ws.Range("A2").Value = 10
$ws.Range("B2").Value = "Vega Modelo de Temuco"
$ws.Range("C2").Value = "La Araucanía"
$ws.Range("D2").Value = 44515
$ws.Range("E2").Value = 9
$ws.Range("F2").Value = "Fruta"
$ws.Range("G2").Value = 100104
$ws.Range("H2").Value = "Frutos de pepita"
$ws.Range("I2").Value = 100104004
$ws.Range("J2").Value = "Níspero"
$ws.Range("K2").Value = "Californiana(o)"
$ws.Range("L2").Value = "Primera"
$ws.Range("M2").Value = 80
$ws.Range("N2").Value = 28000
$ws.Range("O2").Value = 28000
$ws.Range("P2").Value = 28000
$ws.Range("Q2").Value = "$/bandeja 10 kilos"
$ws.Range("R2").Value = "Provincia de Los Andes"
$ws.Range("S2").Value = 2800
$ws.Range("T2").Value = 10
